$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Total Points" (column E) grades for the Customer Class section (rows 3-6)
# matching the max points already recorded in column D.
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Fill in the "Total Points" (column E) grades for the Product Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the active selection to E15 (grading total for Product Class) and scroll back to top
$null = $ws.Range("A1").Select()
$null = $ws.Range("E15").Select()
